# Generate Report for Handback
# Refresh the handback-status report: timestamps move forward as new
# handoff/handback cycles complete, and the zh-cn file's priority flips
# from "ht" (human translation) to "mt" (machine translation).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: bump the "Latest HO Xliff Generate Date" for the
#     01cea992... row (rows 2 and 4 are mirrored entries in this fixture).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-17 22:15:59"
$wsOverview.Range("G4").Value = "2016-08-17 22:15:59"

# --- zh-cn sheet: priority ht -> mt, and refreshed handoff/handback times.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-08-17 22:15:54"
$wsZhCn.Range("H4").Value = "2016-08-17 22:15:54"
$wsZhCn.Range("K2").Value = "2016-08-17 22:16:15"
$wsZhCn.Range("K4").Value = "2016-08-17 22:16:15"

# --- de-de sheet: priority ht -> mt, and refreshed handoff/handback times.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-08-17 22:15:59"
$wsDeDe.Range("H4").Value = "2016-08-17 22:15:59"
$wsDeDe.Range("K2").Value = "2016-08-17 22:16:24"
$wsDeDe.Range("K4").Value = "2016-08-17 22:16:24"
